# Auto-generated edit script: update crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.541.57"
$ws.Range("E2").Value = "  +3.14%  "

$ws.Range("D3").Value = "1.843.90"
$ws.Range("E3").Value = "  +2.40%  "

$ws.Range("D5").Value = "'232.82"
$ws.Range("E5").Value = "  +3.61%  "

$ws.Range("D6").Value = "'0.620"

$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").Value = "'44.54"
$ws.Range("E8").Value = "  +13.71%  "

$ws.Range("D9").Value = "'0.312"
$ws.Range("E9").Value = "  +8.14%  "

$ws.Range("E10").Value = "  +4.46%  "

$ws.Range("E11").Value = "  +2.51%  "

$ws.Range("E12").Value = "  +2.46%  "

$ws.Range("D13").Value = "1.843.16"
$ws.Range("E13").Value = "  +2.42%  "

$ws.Range("D14").Value = "'11.33"
$ws.Range("E14").Value = "  +4.04%  "

$ws.Range("D15").Value = "'0.676"
$ws.Range("E15").Value = "  +7.43%  "

$ws.Range("D16").Value = "'4.73"
$ws.Range("E16").Value = "  +8.46%  "

$ws.Range("D17").Value = "35.519.37"
$ws.Range("E17").Value = "  +3.19%  "

$ws.Range("D18").Value = "'70.46"
$ws.Range("E18").Value = "  +3.57%  "

$ws.Range("E19").Value = "  +5.21%  "

$ws.Range("D20").Value = "'244.30"
$ws.Range("E20").Value = "  +2.21%  "

$ws.Range("D21").Value = "'12.13"
$ws.Range("E21").Value = "  +9.53%  "

$ws.Range("D22").Value = "'4.63"
$ws.Range("E22").Value = "  +13.65%  "

$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").Value = "'2.29"
$ws.Range("E24").Value = "  +6.10%  "

$ws.Range("D25").Value = "'171.75"
$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").Value = "'8.02"
$ws.Range("E26").Value = "  +4.69%  "

$ws.Range("D27").Value = "'17.87"
$ws.Range("E27").Value = "  +1.93%  "

$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").Value = "'1.58"
$ws.Range("E29").Value = "  +29.43%  "

$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("D31").Value = "3.345.42"
$ws.Range("E31").Value = "  +37.69%  "

$ws.Range("E32").Value = "  +8.27%  "

$ws.Range("E33").Value = "  +7.86%  "

$ws.Range("D34").Value = "'3.95"
$ws.Range("E34").Value = "  +5.56%  "

$ws.Range("D35").Value = "'1.85"
$ws.Range("E35").Value = "  +2.37%  "

$ws.Range("D36").Value = "'96.16"
$ws.Range("E36").Value = "  +17.59%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.15"
$ws.Range("E37").Value = "  +9.12%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.694"
$ws.Range("E38").Value = "  +8.52%  "

$ws.Range("D39").Value = "1.346.43"
$ws.Range("E39").Value = "  +2.95%  "

$ws.Range("E40").Value = "  +6.16%  "

$ws.Range("D41").Value = "'2.45"
$ws.Range("E41").Value = "  +6.48%  "

$ws.Range("E42").Value = "  +7.90%  "

$ws.Range("D43").Value = "'15.21"
$ws.Range("E43").Value = "  +8.16%  "

$ws.Range("E44").Value = "  +3.43%  "

$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("D47").Value = "'6.28"
$ws.Range("E47").Value = "  +10.02%  "

$ws.Range("D48").Value = "'0.0518"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("D49").Value = "2.017.17"
$ws.Range("E49").Value = "  +2.90%  "

$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").Value = "'102.51"
$ws.Range("E51").Value = "  +0.67%  "

